$wb = $excel.ActiveWorkbook

# Add the new sheet at the end of the workbook (after the last existing sheet)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "User - Ganti Password"

# Header row
$ws.Range("A1").Value = "currentPassword"
$ws.Range("B1").Value = "newPassword"
$ws.Range("C1").Value = "confirmPassword"
$ws.Range("D1").Value = "otp1"
$ws.Range("E1").Value = "otp2"
$ws.Range("F1").Value = "otp3"
$ws.Range("G1").Value = "otp4"
$ws.Range("H1").Value = "otp5"
$ws.Range("I1").Value = "otp6"
$ws.Range("J1").Value = "kondisi"
$ws.Range("K1").Value = "keterangan"

# Row 2
$ws.Range("J2").Value = "fail"
$ws.Range("K2").Value = "emptyAll"

# Row 3
$ws.Range("B3").Value = "Password3"
$ws.Range("C3").Value = "Password3"
$ws.Range("J3").Value = "fail"
$ws.Range("K3").Value = "emptyCurrentPassword"

# Row 4
$ws.Range("A4").Value = "Password2"
$ws.Range("C4").Value = "Password3"
$ws.Range("J4").Value = "fail"
$ws.Range("K4").Value = "emptyNewPassword"

# Row 5
$ws.Range("A5").Value = "Password2"
$ws.Range("B5").Value = "Password3"
$ws.Range("J5").Value = "fail"
$ws.Range("K5").Value = "emptyConfirmPassword"

# Row 6
$ws.Range("A6").Value = "Password10"
$ws.Range("B6").Value = "Password3"
$ws.Range("C6").Value = "Password3"
$ws.Range("J6").Value = "fail"
$ws.Range("K6").Value = "invalidCurrentPassword"

# Row 7
$ws.Range("A7").Value = "Password2"
$ws.Range("B7").Value = "Pass"
$ws.Range("C7").Value = "Pass"
$ws.Range("J7").Value = "fail"
$ws.Range("K7").Value = "invalidNewPassword"

# Row 8
$ws.Range("A8").Value = "Password2"
$ws.Range("B8").Value = "Password3"
$ws.Range("C8").Value = "Password4"
$ws.Range("J8").Value = "fail"
$ws.Range("K8").Value = "notMatch"

# Row 9
$ws.Range("A9").Value = "Password2"
$ws.Range("B9").Value = "Password3"
$ws.Range("C9").Value = "Password3"
$ws.Range("J9").Value = "fail"
$ws.Range("K9").Value = "cancel"

# Row 10
$ws.Range("A10").Value = "Password2"
$ws.Range("B10").Value = "Password3"
$ws.Range("C10").Value = "Password3"
$ws.Range("D10").Value = 9
$ws.Range("E10").Value = 6
$ws.Range("F10").Value = 8
$ws.Range("G10").Value = 5
$ws.Range("H10").Value = 2
$ws.Range("I10").Value = 6
$ws.Range("J10").Value = "fail"
$ws.Range("K10").Value = "wrongOTP"

# Row 11
$ws.Range("A11").Value = "Password2"
$ws.Range("B11").Value = "Password3"
$ws.Range("C11").Value = "Password3"
$ws.Range("D11").Value = 2
$ws.Range("E11").Value = 4
$ws.Range("F11").Value = 9
$ws.Range("G11").Value = 9
$ws.Range("H11").Value = 2
$ws.Range("I11").Value = 1
$ws.Range("J11").Value = "pass"

# Column widths
$ws.Columns.Item(1).ColumnWidth = 16.109375
$ws.Columns.Item(2).ColumnWidth = 12
$ws.Columns.Item(3).ColumnWidth = 16.21875
$ws.Columns.Item(11).ColumnWidth = 21.88671875

# Selection, matching the saved view state
$ws.Range("K13").Select() | Out-Null
